$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.247.35"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.982.37"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.41"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.28"
$ws.Range("E6").Value = "  -6.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.978.21"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -6.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.93"
$ws.Range("E11").Value = "  -7.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.68"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.461.79"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.253.55"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.110"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.982.05"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.08"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.662"
$ws.Range("E22").Value = "  -5.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.64"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.90"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.34"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.68"
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("E36").Value = "  -3.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "443.86"
$ws.Range("E37").Value = "  -10.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.134.37"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0785"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0377"
$ws.Range("E40").Value = "  -6.62%  "
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  -11.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.47"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("E48").Value = "  -5.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.76"
$ws.Range("E49").Value = "  -7.37%  "
$ws.Range("E50").Value = "  +9.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₃0479"
$ws.Range("E51").Value = "  -10.04%  "
